# Append the "NÖTIGE ÄNDERUNGEN" follow-up notes after the final "Antwort:" paragraph.
$d = $word.ActiveDocument

$antwortIndex = $d.Paragraphs.Count
$antwort = $d.Paragraphs($antwortIndex)

# --- 1) Build the six numbered "Listenabsatz" bullet items right after "Antwort:" ---
# Each clones list/numbering formatting (Listenabsatz, numId 1, ilvl 0) straight
# from the "Antwort:" paragraph, so the existing list (numId=1) is continued
# instead of a brand-new list definition being minted.
$bulletTexts = @(
    "Bei fragetyp wahl => anzeigen wievielte Frage",
    "Dauer der Umfrage angeben",
    "Questions und probably Answers werden nicht in json gepusht",
    "Bei Eingabe antworttitel => anzeigen wievielte antwort",
    "Bei custom questions werden antworten nicht gepusht"
)

$cur = $antwort
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Bei fragetyp wahl => anzeigen wievielte Frage"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Dauer der Umfrage angeben"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Questions und probably Answers werden nicht in json gepusht"
$cur.Range.LanguageID = "en-US"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Bei Eingabe antworttitel => anzeigen wievielte antwort"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Bei custom questions werden antworten nicht gepusht"

# --- 2) Insert the two plain (non-list) paragraphs between "Antwort:" and the ---
#        first bullet: one blank line, then the "NÖTIGE ÄNDERUNGEN:" heading.
$firstBulletIndex = $antwortIndex + 1
$firstBullet = $d.Paragraphs($firstBulletIndex)
$firstBullet.Range.InsertParagraphBefore()

# Fix formatting of the new blank paragraph (strip numbering / revert to
# the document's default "Standard" style) before giving it any text.
$blank = $d.Paragraphs($firstBulletIndex)
$blank.Range.ListFormat.RemoveNumbers()
$blank.Style = $d.Styles("Standard")

# Insert the heading paragraph right after the blank one, before the first bullet.
$firstBullet = $d.Paragraphs($firstBulletIndex + 1)
$firstBullet.Range.InsertParagraphBefore()
$heading = $d.Paragraphs($firstBulletIndex + 1)
$heading.Range.ListFormat.RemoveNumbers()
$heading.Style = $d.Styles("Standard")
$heading.Range.Text = "NÖTIGE ÄNDERUNGEN:"
